# Scheduled-runner refresh of market/profit figures (columns H:N) across the
# per-job-class "Golem Profits" sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ.
# Cells set to "" clear the cell entirely (matches rows where the refreshed
# profit figure is no longer present in the source data).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 382
$ws.Range("I33").Value = 382
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 382
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -153
$ws.Range("N33").Value = ""
$ws.Range("H43").Value = 5903.3335
$ws.Range("I43").Value = 5763
$ws.Range("J43").Value = 6099.8
$ws.Range("K43").Value = 5763
$ws.Range("L43").Value = 6099.8
$ws.Range("M43").Value = -5694
$ws.Range("N43").Value = -6237.8
$ws.Range("H58").Value = 73
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").Value = ""
$ws.Range("H97").Value = 1499.6666
$ws.Range("J97").Value = 1499.6666
$ws.Range("L97").Value = 4498.9998
$ws.Range("N97").Value = -5490.9998
$ws.Range("H132").Value = 2689.3635
$ws.Range("I132").Value = 3028.2222
$ws.Range("J132").Value = 1164.5
$ws.Range("K132").Value = 9084.6666
$ws.Range("L132").Value = 3493.5
$ws.Range("M132").Value = -6554.6666
$ws.Range("N132").Value = -8553.5
$ws.Range("H140").Value = 75000
$ws.Range("J140").Value = 75000
$ws.Range("L140").Value = 75000
$ws.Range("N140").Value = -85360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 5333.3335
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 5333.3335
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 5333.3335
$ws.Range("M21").Value = ""
$ws.Range("N21").Value = -6081.3335
$ws.Range("H30").Value = 13366.5
$ws.Range("I30").Value = 2500
$ws.Range("J30").Value = 15539.8
$ws.Range("K30").Value = 2500
$ws.Range("L30").Value = 15539.8
$ws.Range("M30").Value = -2350
$ws.Range("N30").Value = -15839.8
$ws.Range("H32").Value = 9132.77
$ws.Range("I32").Value = 9132.77
$ws.Range("K32").Value = 9132.77
$ws.Range("M32").Value = -8845.77
$ws.Range("H76").Value = 40244
$ws.Range("J76").Value = 40244
$ws.Range("L76").Value = 40244
$ws.Range("N76").Value = -40920
$ws.Range("H79").Value = 40244
$ws.Range("J79").Value = 40244
$ws.Range("L79").Value = 40244
$ws.Range("N79").Value = -42584
$ws.Range("H96").Value = 33250
$ws.Range("J96").Value = 33250
$ws.Range("L96").Value = 33250
$ws.Range("N96").Value = -38742
$ws.Range("H102").Value = 21001720
$ws.Range("I102").Value = 1113021.1
$ws.Range("J102").Value = 200000000
$ws.Range("K102").Value = 1113021.1
$ws.Range("L102").Value = 200000000
$ws.Range("M102").Value = -1111399.1
$ws.Range("N102").Value = -200003244
$ws.Range("H104").Value = 9999.5
$ws.Range("J104").Value = 9999.5
$ws.Range("L104").Value = 9999.5
$ws.Range("N104").Value = -16987.5
$ws.Range("H113").Value = 35000
$ws.Range("J113").Value = 35000
$ws.Range("L113").Value = 35000
$ws.Range("N113").Value = -43678

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 150
$ws.Range("I11").Value = 150
$ws.Range("K11").Value = 150
$ws.Range("M11").Value = -10
$ws.Range("H16").Value = 1544.5
$ws.Range("I16").Value = 499
$ws.Range("J16").Value = 2590
$ws.Range("K16").Value = 499
$ws.Range("L16").Value = 2590
$ws.Range("M16").Value = -329
$ws.Range("N16").Value = -2930
$ws.Range("H99").Value = 1974.9166
$ws.Range("I99").Value = 1974.9166
$ws.Range("K99").Value = 1974.9166
$ws.Range("M99").Value = -476.9166
$ws.Range("H110").Value = 140000
$ws.Range("J110").Value = 140000
$ws.Range("L110").Value = 140000
$ws.Range("N110").Value = -148180

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 22000
$ws.Range("J28").Value = 22000
$ws.Range("L28").Value = 22000
$ws.Range("N28").Value = -22490
$ws.Range("H31").Value = 6004.3335
$ws.Range("I31").Value = 4006.5
$ws.Range("K31").Value = 4006.5
$ws.Range("M31").Value = -3711.5
$ws.Range("H34").Value = 6004.3335
$ws.Range("I34").Value = 4006.5
$ws.Range("K34").Value = 4006.5
$ws.Range("M34").Value = -3804.5
$ws.Range("H86").Value = 500002500
$ws.Range("I86").Value = 500002500
$ws.Range("K86").Value = 500002500
$ws.Range("M86").Value = -500001377
$ws.Range("H89").Value = 500002500
$ws.Range("I89").Value = 500002500
$ws.Range("K89").Value = 2500012500
$ws.Range("M89").Value = -2500006884
$ws.Range("H99").Value = 1003999.8
$ws.Range("I99").Value = 4833
$ws.Range("J99").Value = 2502750
$ws.Range("K99").Value = 4833
$ws.Range("L99").Value = 2502750
$ws.Range("M99").Value = -3335
$ws.Range("N99").Value = -2505746
$ws.Range("H126").Value = 1003999.8
$ws.Range("I126").Value = 4833
$ws.Range("J126").Value = 2502750
$ws.Range("K126").Value = 14499
$ws.Range("L126").Value = 7508250
$ws.Range("M126").Value = -12029
$ws.Range("N126").Value = -7513190
$ws.Range("H141").Value = 695183.5600000001
$ws.Range("J141").Value = 695183.5600000001
$ws.Range("L141").Value = 695183.5600000001
$ws.Range("N141").Value = -705543.5600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1541
$ws.Range("I4").Value = 1191.24
$ws.Range("K4").Value = 3573.72
$ws.Range("M4").Value = -3461.72
$ws.Range("H12").Value = 88.166664
$ws.Range("I12").Value = 9
$ws.Range("K12").Value = 27
$ws.Range("M12").Value = 146
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = ""
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = ""
$ws.Range("H99").Value = 8510.286
$ws.Range("I99").Value = 6914.6
$ws.Range("K99").Value = 20743.8
$ws.Range("M99").Value = -18497.8
$ws.Range("H103").Value = 228.42857
$ws.Range("I103").Value = 183.16667
$ws.Range("J103").Value = 500
$ws.Range("K103").Value = 549.50001
$ws.Range("L103").Value = 1500
$ws.Range("M103").Value = 329.49999
$ws.Range("N103").Value = -3258
$ws.Range("H113").Value = 700
$ws.Range("J113").Value = 700
$ws.Range("L113").Value = 2100
$ws.Range("N113").Value = -6440

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").Value = ""
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").Value = ""
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").Value = ""
$ws.Range("H95").Value = 31000
$ws.Range("J95").Value = 31000
$ws.Range("L95").Value = 31000
$ws.Range("N95").Value = -36492
$ws.Range("H98").Value = 42648.6
$ws.Range("J98").Value = 42648.6
$ws.Range("L98").Value = 42648.6
$ws.Range("N98").Value = -48638.6
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = ""
$ws.Range("N126").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 889
$ws.Range("I7").Value = 889
$ws.Range("K7").Value = 889
$ws.Range("M7").Value = -777
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = ""
$ws.Range("N22").Value = ""
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = ""
$ws.Range("N27").Value = ""
$ws.Range("H40").Value = 2823.5
$ws.Range("J40").Value = 3250
$ws.Range("L40").Value = 3250
$ws.Range("N40").Value = -3522
$ws.Range("H82").Value = 6223.75
$ws.Range("I82").Value = 4997.5
$ws.Range("J82").Value = 7450
$ws.Range("K82").Value = 4997.5
$ws.Range("L82").Value = 7450
$ws.Range("M82").Value = -4636.5
$ws.Range("N82").Value = -8172
$ws.Range("H85").Value = 6223.75
$ws.Range("I85").Value = 4997.5
$ws.Range("J85").Value = 7450
$ws.Range("K85").Value = 4997.5
$ws.Range("L85").Value = 7450
$ws.Range("M85").Value = -3749.5
$ws.Range("N85").Value = -9946
$ws.Range("H94").Value = 42000
$ws.Range("J94").Value = 42000
$ws.Range("L94").Value = 42000
$ws.Range("N94").Value = -43352
$ws.Range("H106").Value = 35624
$ws.Range("J106").Value = 35624
$ws.Range("L106").Value = 35624
$ws.Range("N106").Value = -38148
$ws.Range("H126").Value = 889
$ws.Range("I126").Value = 889
$ws.Range("K126").Value = 2667
$ws.Range("M126").Value = -197

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 33749.75
$ws.Range("J69").Value = 33749.75
$ws.Range("L69").Value = 33749.75
$ws.Range("N69").Value = -35247.75
$ws.Range("H72").Value = 33749.75
$ws.Range("J72").Value = 33749.75
$ws.Range("L72").Value = 101249.25
$ws.Range("N72").Value = -108737.25
$ws.Range("H126").Value = 4720.4
$ws.Range("I126").Value = 1199
$ws.Range("K126").Value = 3597
$ws.Range("M126").Value = -1127
$ws.Range("H136").Value = 2699.5
$ws.Range("I136").Value = 2614
$ws.Range("J136").Value = 2956
$ws.Range("K136").Value = 7842
$ws.Range("L136").Value = 8868
$ws.Range("M136").Value = -5292
$ws.Range("N136").Value = -13968
